$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Replace the 1000-sample junction-flooding data block (rows 2-5)
#    with the new reading, and drop the now-unused 5th data row
#    (old row 6) so the sheet shrinks from A1:AH6 to A1:AH5.
# ------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

$newData = @(
    @(45072.50694444445, 18.737, 12.527, 4.029, 39.827, 31.89, 14.745, 46.499, 22.687, 9.44, 14.208, 15.679, 16.314, 4.706, 14.663, 20.445, 12.568, 3.441, 2.244, 215.823, 40.754, 13.534, 26.795, 13.632, 3.069, 23.634, 11.955, 10.86, 12.744, 16.181, 3.455, 41.233, 7.422, 16.92),
    @(45072.51388888889, 14.893, 10.472, 1.759, 32.18, 25.957, 11.721, 45.586, 18.033, 7.766, 11.436, 12.868, 13.434, 3.745, 11.655, 16.406, 10.076, 1.536, 1.025, 170.083, 32.635, 10.758, 21.603, 11.2, 2.132, 22.128, 9.502000000000001, 8.615, 10.099, 13.352, 1.265, 41.503, 5.93, 13.45),
    @(45072.52083333334, 22.1, 16.138, 1.512, 47.976, 39.183, 17.391, 66.31399999999999, 26.759, 11.803, 17.46, 19.247, 20.238, 5.555, 17.294, 24.548, 14.656, 1.096, 0.989, 255.923, 48.346, 15.963, 32.394, 16.946, 2.68, 32.372, 14.1, 12.567, 14.762, 20.126, 0.773, 60.151, 8.948, 19.958),
    @(45072.52777777778, 23.54, 17.34, 1.34, 51.15, 41.91, 18.53, 72.17, 28.5, 12.65, 18.73, 20.53, 21.63, 5.92, 18.42, 26.21, 15.53, 0.87, 0.93, 273.1, 51.53, 17, 34.63, 18.14, 2.73, 35.03, 15.02, 13.33, 15.66, 21.51, 0.5600000000000001, 65.51000000000001, 9.57, 21.26)
)

for ($r = 0; $r -lt $newData.Length; $r++) {
    $rowValues = $newData[$r]
    $targetRow = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $rowValues[$c]
    }
}

# ------------------------------------------------------------------
# 2) Widen several data columns by one character (custom accuracy
#    formatting needs the extra digit of room).
#    ColumnWidth is Excel's character-width units; the stored OOXML
#    width is ColumnWidth + 5/6, so subtract that offset to land on
#    the intended whole-number stored width.
# ------------------------------------------------------------------
$wideCols8 = @(2, 3, 6, 7, 9, 10, 11, 12, 13, 15, 16, 17, 22, 24, 27, 28, 29, 30, 34)
foreach ($colIdx in $wideCols8) {
    $ws.Columns.Item($colIdx).ColumnWidth = 7.166666666666667
}
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
